# Recalculated market-board profit figures for several Leve rows across
# the ALC/ARM/BSM/CUL/GSM/LTW/WVR sheets (currentAveragePrice* / LevePrice* /
# LeveProfit* columns H:N), refreshed by the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6: Days of Chunder (Antidote)
$ws.Range("H6").Value = 134.6923
$ws.Range("I6").Value = 126.2
$ws.Range("J6").Value = 163
$ws.Range("K6").Value = 378.6
$ws.Range("L6").Value = 489
$ws.Range("M6").Value = -266.6
$ws.Range("N6").Value = -713
# Row 11: Gotta Bounce (Rubber)
$ws.Range("H11").Value = 63.4
$ws.Range("I11").Value = 63.4
$ws.Range("K11").Value = 63.4
$ws.Range("M11").Value = 76.59999999999999
# Row 12: Don't Be So Tallow (Beeswax)
$ws.Range("H12").Value = 1000
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 1000
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 1000
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -1340
# Row 41: The Write Stuff (Enchanted Mythril Ink)
$ws.Range("H41").Value = 2825
$ws.Range("I41").Value = 600
$ws.Range("J41").Value = 3142.8572
$ws.Range("K41").Value = 600
$ws.Range("L41").Value = 3142.8572
$ws.Range("M41").Value = -160
$ws.Range("N41").Value = -4022.8572
# Row 112: Making Ends Meet (Superior Spiritbond Potion)
$ws.Range("H112").Value = 2837.5
$ws.Range("J112").Value = 2885.7144
$ws.Range("L112").Value = 8657.143199999999
$ws.Range("N112").Value = -10873.1432
# Row 137: Cutting Edge of Culinary Quality (Magnesia Whetstone)
$ws.Range("H137").Value = 3456.3572
$ws.Range("I137").Value = 2232
$ws.Range("J137").Value = 4374.625
$ws.Range("K137").Value = 6696
$ws.Range("L137").Value = 13123.875
$ws.Range("M137").Value = -4146
$ws.Range("N137").Value = -18223.875

$ws = $wb.Worksheets.Item("ARM")
# Row 45: Hollow Hallmarks (Mythril Ingot)
$ws.Range("H45").Value = 2918.0908
$ws.Range("I45").Value = 2327.8333
$ws.Range("J45").Value = 3626.4
$ws.Range("K45").Value = 2327.8333
$ws.Range("L45").Value = 3626.4
$ws.Range("M45").Value = -1950.8333
$ws.Range("N45").Value = -4380.4
# Row 61: Dealing with the Tough Stuff (Cobalt Ingot)
$ws.Range("H61").Value = 3719.0667
$ws.Range("J61").Value = 6599.8
$ws.Range("L61").Value = 6599.8
$ws.Range("N61").Value = -7023.8
# Row 132: Don't Bore Me, Ore Me (Mountain Chromite Ingot)
$ws.Range("H132").Value = 3333.2778
$ws.Range("I132").Value = 3333.2778
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 9999.8334
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7469.8334
$ws.Range("N132").ClearContents()
# Row 136: Metal with Mettle (Cobalt Tungsten Ingot)
$ws.Range("H136").Value = 3719.0667
$ws.Range("J136").Value = 6599.8
$ws.Range("L136").Value = 19799.4
$ws.Range("N136").Value = -24899.4

$ws = $wb.Worksheets.Item("BSM")
# Row 105: Ingot to Wing It (Molybdenum Ingot)
$ws.Range("H105").Value = 2573.238
$ws.Range("J105").Value = 2791.3333
$ws.Range("L105").Value = 2791.3333
$ws.Range("N105").Value = -6285.3333
# Row 107: The Gold Experience (Deepgold Nugget)
$ws.Range("H107").Value = 5047.0713
$ws.Range("I107").Value = 4338.5713
$ws.Range("K107").Value = 4338.5713
$ws.Range("M107").Value = -2418.5713
# Row 134: Ruthenium Supremium (Ruthenium Ingot)
$ws.Range("H134").Value = 1419.8334
$ws.Range("I134").Value = 1419.8334
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 4259.5002
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -1724.5002
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 131: The Mountain Steeped (Tsai tou Vounou)
$ws.Range("H131").Value = 1645.5
$ws.Range("J131").Value = 2315
$ws.Range("L131").Value = 6945
$ws.Range("N131").Value = -17025
# Row 138: Bring Me Your Tacos (Tacos Al Pastor)
$ws.Range("H138").Value = 4498.875
$ws.Range("I138").Value = 2998.5
$ws.Range("J138").Value = 9000
$ws.Range("K138").Value = 8995.5
$ws.Range("L138").Value = 27000
$ws.Range("M138").Value = -3855.5
$ws.Range("N138").Value = -37280
# Row 140: Sweet, Sweet Bean Juice (Mesquite Juice)
$ws.Range("H140").Value = 2947.6365
$ws.Range("I140").Value = 2742.4
$ws.Range("J140").Value = 5000
$ws.Range("K140").Value = 8227.200000000001
$ws.Range("L140").Value = 15000
$ws.Range("M140").Value = -3047.200000000001
$ws.Range("N140").Value = -25360

$ws = $wb.Worksheets.Item("GSM")
# Row 31: One and Only (Staghorn Staff)
$ws.Range("H31").Value = 2544
$ws.Range("I31").Value = 2107.625
$ws.Range("K31").Value = 2107.625
$ws.Range("M31").Value = -1815.625
# Row 37: Dancing with the Stars (Toothed Staghorn Staff)
$ws.Range("H37").Value = 2544
$ws.Range("I37").Value = 2107.625
$ws.Range("K37").Value = 2107.625
$ws.Range("M37").Value = -1830.625
# Row 80: Needs More Prayerbell (Hardsilver Ingot)
$ws.Range("H80").Value = 1166.6666
$ws.Range("I80").Value = 1166.6666
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 1166.6666
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -168.6666
$ws.Range("N80").ClearContents()
# Row 83: With a Noise That Reaches Heaven (L) (Hardsilver Ingot)
$ws.Range("H83").Value = 1166.6666
$ws.Range("I83").Value = 1166.6666
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 5833.333000000001
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -841.3330000000005
$ws.Range("N83").ClearContents()
# Row 132: On Board for Lar (Lar Ingot)
$ws.Range("H132").Value = 3501.5356
$ws.Range("I132").Value = 2686.7827
$ws.Range("J132").Value = 7249.4
$ws.Range("K132").Value = 8060.348100000001
$ws.Range("L132").Value = 21748.2
$ws.Range("M132").Value = -5530.348100000001
$ws.Range("N132").Value = -26808.2

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban (Leather)
$ws.Range("H7").Value = 5952.7646
$ws.Range("I7").Value = 4944
$ws.Range("K7").Value = 4944
$ws.Range("M7").Value = -4832
# Row 16: Saddle Sore (Hard Leather)
$ws.Range("H16").Value = 671.4286
$ws.Range("I16").Value = 725.25
$ws.Range("K16").Value = 725.25
$ws.Range("M16").Value = -555.25
# Row 55: It's Not a Job, It's a Calling (Peiste Leather)
$ws.Range("H55").Value = 1541.2858
$ws.Range("I55").Value = 1541.2858
$ws.Range("K55").Value = 1541.2858
$ws.Range("M55").Value = -1368.2858
# Row 68: You Could Say It's a Moving Target (Wyvern Leather)
$ws.Range("H68").Value = 6443.3335
$ws.Range("I68").Value = 2497.5
$ws.Range("K68").Value = 2497.5
$ws.Range("M68").Value = -1748.5
# Row 71: They Call It Bloody Mary (L) (Wyvern Leather)
$ws.Range("H71").Value = 6443.3335
$ws.Range("I71").Value = 2497.5
$ws.Range("K71").Value = 12487.5
$ws.Range("M71").Value = -8743.5
# Row 74: Overall, We Blend In (Dhalmelskin Vest)
$ws.Range("H74").Value = 18750
$ws.Range("I74").Value = 18750
$ws.Range("K74").Value = 18750
$ws.Range("M74").Value = -17752
# Row 77: Eviction Notice (L) (Dhalmelskin Vest)
$ws.Range("H77").Value = 18750
$ws.Range("I77").Value = 18750
$ws.Range("K77").Value = 56250
$ws.Range("M77").Value = -51258
# Row 126: Battered Books (Saiga Leather)
$ws.Range("H126").Value = 5952.7646
$ws.Range("I126").Value = 4944
$ws.Range("K126").Value = 14832
$ws.Range("M126").Value = -12362

$ws = $wb.Worksheets.Item("WVR")
# Row 54: No Country for Cold Men (Woolen Tights)
$ws.Range("H54").Value = 59999.5
$ws.Range("J54").Value = 59999.5
$ws.Range("L54").Value = 59999.5
$ws.Range("N54").Value = -61039.5
# Row 62: Pride Up in Smoke (Rainbow Cloth)
$ws.Range("H62").Value = 12230.923
$ws.Range("J62").Value = 12818.182
$ws.Range("L62").Value = 12818.182
$ws.Range("N62").Value = -14066.182
# Row 65: Desperate for Diversionaries (L) (Rainbow Cloth)
$ws.Range("H65").Value = 12230.923
$ws.Range("J65").Value = 12818.182
$ws.Range("L65").Value = 64090.91
$ws.Range("N65").Value = -70330.91
# Row 136: Weaving the Envelope (Sarcenet Cloth)
$ws.Range("H136").Value = 3892.4285
$ws.Range("I136").Value = 2943.889
$ws.Range("J136").Value = 5599.8
$ws.Range("K136").Value = 8831.667000000001
$ws.Range("L136").Value = 16799.4
$ws.Range("M136").Value = -6281.667000000001
$ws.Range("N136").Value = -21899.4
